# Renumber NPC skill rows 11-37 on the "Skill" sheet from the "13_" series to
# the "14_" series (cid + name), keeping style/enhancer/cooler/relic tube
# references untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Skill")

$names = @(
    "14_cat_01_00",
    "14_cat_01_01",
    "14_cat_02_00",
    "14_cat_02_00",
    "14_pitbull_00",
    "14_pitbull_01",
    "14_pitbull_02",
    "14_owl_00",
    "14_owl_01",
    "14_crow_00",
    "14_crow_01",
    "14_salamander_00",
    "14_carbannog_00",
    "14_carbannog_01",
    "14_carbannog_02",
    "14_carbannog_03",
    "14_tag_00",
    "14_tag_01",
    "14_tag_02",
    "14_tim_00",
    "14_tim_01",
    "14_tim_02",
    "14_tim_03",
    "14_madbuddy_00",
    "14_madbuddy_01",
    "14_madbuddy_02",
    "14_madbuddy_03"
)

$startRow = 11
$cid = 14100

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $cid
    $ws.Cells.Item($row, 2).Value = $names[$i]
    $cid = $cid + 1
}

$ws.Activate()
$ws.Range("B16").Select()
